$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (row 2 and row 3) need to have their values swapped
# for columns D, M, N, O, P, Q, S, T, turning the dataset into
# chronological (weekly) order.

$cols = @("D", "M", "N", "O", "P", "Q", "S", "T")

foreach ($col in $cols) {
    $cellA = $ws.Range($col + "2")
    $cellB = $ws.Range($col + "3")
    $tmp = $cellA.Value2
    $cellA.Value2 = $cellB.Value2
    $cellB.Value2 = $tmp
}
